# "Changes in code for scope 1"
# Update the Natural Gas row's CH4 and N2O emission factors, then leave the
# selection where the author left off (cell E11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Natural Gas row (row 4): CH4 Factor (D4) and N2O Factor (E4)
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.1

# Match the workbook's final selection state
$ws.Activate()
$ws.Range("E11").Select()
